$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 336; existing rows 336-374 shift down to 337-375.
$ws.Rows.Item(336).Insert()

# Populate the newly inserted row 336 with its data.
$ws.Cells.Item(336, 1).Value = 3
$ws.Cells.Item(336, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(336, 3).Value = "Coquimbo"
$ws.Cells.Item(336, 4).Value = Get-Date -Year 2023 -Month 9 -Day 25 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(336, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(336, 5).Value = 5
$ws.Cells.Item(336, 6).Value = "Fruta"
$ws.Cells.Item(336, 7).Value = 100101
$ws.Cells.Item(336, 8).Value = "Berries"
$ws.Cells.Item(336, 9).Value = 100101001
$ws.Cells.Item(336, 10).Value = "Arándano (blue)"
$ws.Cells.Item(336, 11).Value = "Sin especificar"
$ws.Cells.Item(336, 12).Value = "Primera"
$ws.Cells.Item(336, 13).Value = 36
$ws.Cells.Item(336, 14).Value = 16000
$ws.Cells.Item(336, 15).Value = 16000
$ws.Cells.Item(336, 16).Value = 16000
$ws.Cells.Item(336, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(336, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(336, 19).Value = 8000
$ws.Cells.Item(336, 20).Value = 2
